# Add 4 new course / course-package rows to the "courses" sheet
# (RII50520, RII60520, ICT60220 and a civil-construction-design package).
#
# NOTE: cell values are written in a specific sequence on purpose so that
# the workbook's shared-string table and style table end up populated in
# the same order as the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: RII50520 - Diploma of Civil Construction Design ----
$ws.Range("H2").WrapText = $true
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("M2").Value = "TAS"
$ws.Range("R2").Style = "Normal"
$ws.Range("R2").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("A2").Value = "RII50520"
$ws.Range("B2").Value = "111827M"
$ws.Range("I2").NumberFormat = "#,##0"
$ws.Range("I2").Value = 10200
$ws.Range("J2").WrapText = $true
$ws.Range("J2").NumberFormat = "#,##0"
$ws.Range("J2").Value = "10,000 tuition fee + 200 handling fee"

# ---- Row 3: RII60520 - Advanced Diploma of Civil Construction Design ----
$ws.Range("D3").Value = "ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("D2").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("B3").Value = "111826A"
$ws.Range("A3").Value = "RII60520"
$ws.Range("H3").WrapText = $true
$ws.Range("H3").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("J3").WrapText = $true
$ws.Range("J3").NumberFormat = "#,##0"
$ws.Range("J3").Value = "20,000 tuition fee + 200 handling fee"

# ---- Row 4: ICT60220 - Advanced Diploma of Information Technology ----
$ws.Range("D4").WrapText = $true
$ws.Range("D4").Value = "ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY (Telecommunications Network Engineering)"
$ws.Range("A4").Value = "ICT60220"
$ws.Range("B4").Value = "111825B"
$ws.Range("J4").WrapText = $true
$ws.Range("J4").NumberFormat = "#,##0"
$ws.Range("J4").Value = "16,000 tuition fee + 200 handling fee"

# ---- Row 5: Package - Diploma + Advanced Diploma of Civil Construction Design ----
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("D5").WrapText = $true
$ws.Range("D5").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Value = "RII50520/RII60520"
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Value = "111827M/111826A"
$ws.Range("C2").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("C4").Value = "INFORMATION TECHNOLOGY"

# ---- Remaining numeric cells and cells that reuse existing shared strings ----
$ws.Range("E2").Value = 52
$ws.Range("E3").Value = 104
$ws.Range("I3").NumberFormat = "#,##0"
$ws.Range("I3").Value = 20200
$ws.Range("E4").Value = 104
$ws.Range("I4").NumberFormat = "#,##0"
$ws.Range("I4").Value = 16200
$ws.Range("E5").Value = 104
$ws.Range("I5").NumberFormat = "#,##0"
$ws.Range("I5").Value = 20200
$ws.Range("C3").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("H4").WrapText = $true
$ws.Range("H4").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("H5").WrapText = $true
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("J5").WrapText = $true
$ws.Range("J5").NumberFormat = "#,##0"
$ws.Range("J5").Value = "20,000 tuition fee + 200 handling fee"
$ws.Range("M3").Value = "TAS"
$ws.Range("M4").Value = "TAS"
$ws.Range("M5").Value = "TAS"
$ws.Range("R3").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R4").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "
$ws.Range("R5").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# ---- Row heights for the new rows ----
$ws.Rows(2).RowHeight = 45
$ws.Rows(3).RowHeight = 45
$ws.Rows(4).RowHeight = 45
$ws.Rows(5).RowHeight = 45

# ---- Restore the active selection as left by the author ----
$ws.Range("E8").Select()
